$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update data rows 2-5 with the new dataset (custom accuracy + 1000 data points) ----
$row2 = New-Object "object[,]" 1,34
$row2[0,0] = 45085.50694444445
$row2[0,1] = 15.374
$row2[0,2] = 10.112
$row2[0,3] = 3.717
$row2[0,4] = 32.732
$row2[0,5] = 25.126
$row2[0,6] = 11.943
$row2[0,7] = 36.363
$row2[0,8] = 18.615
$row2[0,9] = 7.558
$row2[0,10] = 11.183
$row2[0,11] = 12.932
$row2[0,12] = 13.607
$row2[0,13] = 3.86
$row2[0,14] = 12.031
$row2[0,15] = 16.608
$row2[0,16] = 10.57
$row2[0,17] = 3.118
$row2[0,18] = 1.764
$row2[0,19] = 175.75
$row2[0,20] = 33.355
$row2[0,21] = 11.105
$row2[0,22] = 21.604
$row2[0,23] = 11.537
$row2[0,24] = 2.926
$row2[0,25] = 18.421
$row2[0,26] = 9.808999999999999
$row2[0,27] = 8.874000000000001
$row2[0,28] = 10.611
$row2[0,29] = 13.63
$row2[0,30] = 3.312
$row2[0,31] = 32.661
$row2[0,32] = 5.894
$row2[0,33] = 13.883
$ws.Range("A2:AH2").Value = $row2

$row3 = New-Object "object[,]" 1,34
$row3[0,0] = 45085.51388888889
$row3[0,1] = 2.402
$row3[0,2] = 1.114
$row3[0,3] = 1.307
$row3[0,4] = 5.089
$row3[0,5] = 3.122
$row3[0,6] = 1.767
$row3[0,7] = 11.949
$row3[0,8] = 2.909
$row3[0,9] = 1.058
$row3[0,10] = 1.171
$row3[0,11] = 1.983
$row3[0,12] = 2.086
$row3[0,13] = 0.625
$row3[0,14] = 1.88
$row3[0,15] = 2.541
$row3[0,16] = 2.039
$row3[0,17] = 1.354
$row3[0,18] = 0.609
$row3[0,19] = 21.385
$row3[0,20] = 5.67
$row3[0,21] = 1.735
$row3[0,22] = 3.405
$row3[0,23] = 1.867
$row3[0,24] = 0.889
$row3[0,25] = 5.321
$row3[0,26] = 1.533
$row3[0,27] = 1.584
$row3[0,28] = 1.889
$row3[0,29] = 2.087
$row3[0,30] = 1.246
$row3[0,31] = 11.719
$row3[0,32] = 0.725
$row3[0,33] = 2.176
$ws.Range("A3:AH3").Value = $row3

$row4 = New-Object "object[,]" 1,34
$row4[0,0] = 45085.52083333334
$row4[0,1] = 18.737
$row4[0,2] = 13.622
$row4[0,3] = 1.386
$row4[0,4] = 40.703
$row4[0,5] = 32.891
$row4[0,6] = 14.639
$row4[0,7] = 52.254
$row4[0,8] = 22.687
$row4[0,9] = 9.951000000000001
$row4[0,10] = 14.666
$row4[0,11] = 16.311
$row4[0,12] = 17.273
$row4[0,13] = 4.708
$row4[0,14] = 14.663
$row4[0,15] = 20.746
$row4[0,16] = 12.505
$row4[0,17] = 1.028
$row4[0,18] = 0.878
$row4[0,19] = 215.837
$row4[0,20] = 40.861
$row4[0,21] = 13.534
$row4[0,22] = 27.29
$row4[0,23] = 14.523
$row4[0,24] = 2.351
$row4[0,25] = 25.83
$row4[0,26] = 11.955
$row4[0,27] = 10.661
$row4[0,28] = 12.555
$row4[0,29] = 17.144
$row4[0,30] = 0.766
$row4[0,31] = 47.355
$row4[0,32] = 7.53
$row4[0,33] = 16.92
$ws.Range("A4:AH4").Value = $row4

$row5 = New-Object "object[,]" 1,34
$row5[0,0] = 45085.52777777778
$row5[0,1] = 12.01
$row5[0,2] = 8.69
$row5[0,3] = 0.95
$row5[0,4] = 26.1
$row5[0,5] = 20.99
$row5[0,6] = 9.359999999999999
$row5[0,7] = 38.45
$row5[0,8] = 14.54
$row5[0,9] = 6.38
$row5[0,10] = 9.32
$row5[0,11] = 10.46
$row5[0,12] = 11.1
$row5[0,13] = 3.02
$row5[0,14] = 9.4
$row5[0,15] = 13.31
$row5[0,16] = 8.07
$row5[0,17] = 0.75
$row5[0,18] = 0.58
$row5[0,19] = 135.75
$row5[0,20] = 26.36
$row5[0,21] = 8.68
$row5[0,22] = 17.56
$row5[0,23] = 9.33
$row5[0,24] = 1.55
$row5[0,25] = 18.33
$row5[0,26] = 7.66
$row5[0,27] = 6.86
$row5[0,28] = 8.07
$row5[0,29] = 10.99
$row5[0,30] = 0.55
$row5[0,31] = 35.21
$row5[0,32] = 4.8
$row5[0,33] = 10.85
$ws.Range("A5:AH5").Value = $row5


# ---- Widen several data columns from 7 to 8 characters ----
$ws.Columns.Item(2).ColumnWidth = 7.17
$ws.Columns.Item(3).ColumnWidth = 7.17
$ws.Columns.Item(5).ColumnWidth = 7.17
$ws.Columns.Item(7).ColumnWidth = 7.17
$ws.Columns.Item(11).ColumnWidth = 7.17
$ws.Columns.Item(12).ColumnWidth = 7.17
$ws.Columns.Item(13).ColumnWidth = 7.17
$ws.Columns.Item(15).ColumnWidth = 7.17
$ws.Columns.Item(17).ColumnWidth = 7.17
$ws.Columns.Item(21).ColumnWidth = 7.17
$ws.Columns.Item(22).ColumnWidth = 7.17
$ws.Columns.Item(24).ColumnWidth = 7.17
$ws.Columns.Item(27).ColumnWidth = 7.17
$ws.Columns.Item(28).ColumnWidth = 7.17
$ws.Columns.Item(29).ColumnWidth = 7.17
$ws.Columns.Item(30).ColumnWidth = 7.17
$ws.Columns.Item(34).ColumnWidth = 7.17


# ---- Remove the now-obsolete trailing row (row 6) ----
$ws.Rows.Item(6).Delete()
